$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col8a1"
$ws.Range("C2").Value = "Itga1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.430265
$ws.Range("H2").Value = 34.290795
$ws.Range("I2").Value = 0.05939319992829472
$ws.Range("J2").Value = 0.05939319992829471
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 43.19793166666667
$ws.Range("N2").Value = 129.593795
$ws.Range("O2").Value = 0.7412538312889448
$ws.Range("P2").Value = 0.7412538312889448
$ws.Range("Q2").Value = 493.7638064018917
$ws.Range("R2").Value = 4443.874257617024
$ws.Range("S2").Value = 0.04402543699935874
$ws.Range("T2").Value = 0.04402543699935874

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col8a1"
$ws.Range("C3").Value = "Itga1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.430265
$ws.Range("H3").Value = 34.290795
$ws.Range("I3").Value = 0.05939319992829472
$ws.Range("J3").Value = 0.05939319992829471
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.38032666666667
$ws.Range("N3").Value = 31.14098
$ws.Range("O3").Value = 0.1781209566020688
$ws.Range("P3").Value = 0.1781209566020688
$ws.Range("Q3").Value = 118.6498845865667
$ws.Range("R3").Value = 1067.8489612791
$ws.Range("S3").Value = 0.01057917358688578
$ws.Range("T3").Value = 0.01057917358688578

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col8a1"
$ws.Range("C4").Value = "Itga1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.430265
$ws.Range("H4").Value = 34.290795
$ws.Range("I4").Value = 0.05939319992829472
$ws.Range("J4").Value = 0.05939319992829471
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.422089
$ws.Range("N4").Value = 1.266267
$ws.Range("O4").Value = 0.007242825670663926
$ws.Range("P4").Value = 0.007242825670663927
$ws.Range("Q4").Value = 4.824589123585
$ws.Range("R4").Value = 43.421302112265
$ws.Range("S4").Value = 0.0004301745931035278
$ws.Range("T4").Value = 0.0004301745931035278

# Row 5: ECs -> MuSCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col8a1"
$ws.Range("C5").Value = "Itga1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.430265
$ws.Range("H5").Value = 34.290795
$ws.Range("I5").Value = 0.05939319992829472
$ws.Range("J5").Value = 0.05939319992829471
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.133026333333333
$ws.Range("N5").Value = 12.399079
$ws.Range("O5").Value = 0.07092056230936286
$ws.Range("P5").Value = 0.07092056230936288
$ws.Range("Q5").Value = 47.24158624197832
$ws.Range("R5").Value = 425.1742761778049
$ws.Range("S5").Value = 0.004212199136267071
$ws.Range("T5").Value = 0.004212199136267072

# Row 6: ECs -> Resolving-Mac
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Col8a1"
$ws.Range("C6").Value = "Itga1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.430265
$ws.Range("H6").Value = 34.290795
$ws.Range("I6").Value = 0.05939319992829472
$ws.Range("J6").Value = 0.05939319992829471
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1434673333333333
$ws.Range("N6").Value = 0.430402
$ws.Range("O6").Value = 0.002461824128959449
$ws.Range("P6").Value = 0.002461824128959449
$ws.Range("Q6").Value = 1.639869638843333
$ws.Range("R6").Value = 14.75882674959
$ws.Range("S6").Value = 0.0001462156126795886
$ws.Range("T6").Value = 0.0001462156126795886

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col8a1"
$ws.Range("C7").Value = "Itga1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 175.4151816666667
$ws.Range("H7").Value = 526.245545
$ws.Range("I7").Value = 0.9114809634935387
$ws.Range("J7").Value = 0.9114809634935386
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 43.19793166666667
$ws.Range("N7").Value = 129.593795
$ws.Range("O7").Value = 0.7412538312889448
$ws.Range("P7").Value = 0.7412538312889448
$ws.Range("Q7").Value = 7577.573030932586
$ws.Range("R7").Value = 68198.15727839328
$ws.Range("S7").Value = 0.6756387563365244
$ws.Range("T7").Value = 0.6756387563365244

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col8a1"
$ws.Range("C8").Value = "Itga1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 175.4151816666667
$ws.Range("H8").Value = 526.245545
$ws.Range("I8").Value = 0.9114809634935387
$ws.Range("J8").Value = 0.9114809634935386
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.38032666666667
$ws.Range("N8").Value = 31.14098
$ws.Range("O8").Value = 0.1781209566020688
$ws.Range("P8").Value = 0.1781209566020688
$ws.Range("Q8").Value = 1820.866887992678
$ws.Range("R8").Value = 16387.8019919341
$ws.Range("S8").Value = 0.1623538611420445
$ws.Range("T8").Value = 0.1623538611420445

# Row 9: FAPs -> Inflammatory-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col8a1"
$ws.Range("C9").Value = "Itga1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 175.4151816666667
$ws.Range("H9").Value = 526.245545
$ws.Range("I9").Value = 0.9114809634935387
$ws.Range("J9").Value = 0.9114809634935386
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.422089
$ws.Range("N9").Value = 1.266267
$ws.Range("O9").Value = 0.007242825670663926
$ws.Range("P9").Value = 0.007242825670663927
$ws.Range("Q9").Value = 74.04081861450166
$ws.Range("R9").Value = 666.367367530515
$ws.Range("S9").Value = 0.006601697720712491
$ws.Range("T9").Value = 0.006601697720712491

# Row 10: FAPs -> MuSCs
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Col8a1"
$ws.Range("C10").Value = "Itga1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 175.4151816666667
$ws.Range("H10").Value = 526.245545
$ws.Range("I10").Value = 0.9114809634935387
$ws.Range("J10").Value = 0.9114809634935386
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.133026333333333
$ws.Range("N10").Value = 12.399079
$ws.Range("O10").Value = 0.07092056230936286
$ws.Range("P10").Value = 0.07092056230936288
$ws.Range("Q10").Value = 724.9955650947837
$ws.Range("R10").Value = 6524.960085853054
$ws.Range("S10").Value = 0.06464274246524161
$ws.Range("T10").Value = 0.06464274246524161

# Row 11: FAPs -> Resolving-Mac
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Col8a1"
$ws.Range("C11").Value = "Itga1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 175.4151816666667
$ws.Range("H11").Value = 526.245545
$ws.Range("I11").Value = 0.9114809634935387
$ws.Range("J11").Value = 0.9114809634935386
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.1434673333333333
$ws.Range("N11").Value = 0.430402
$ws.Range("O11").Value = 0.002461824128959449
$ws.Range("P11").Value = 0.002461824128959449
$ws.Range("Q11").Value = 25.16634833989889
$ws.Range("R11").Value = 226.49713505909
$ws.Range("S11").Value = 0.002243905829015601
$ws.Range("T11").Value = 0.002243905829015601

# Row 12: Inflammatory-Mac -> ECs
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Col8a1"
$ws.Range("C12").Value = "Itga1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.05623966666666667
$ws.Range("H12").Value = 0.168719
$ws.Range("I12").Value = 0.0002922288998753735
$ws.Range("J12").Value = 0.0002922288998753735
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 43.19793166666667
$ws.Range("N12").Value = 129.593795
$ws.Range("O12").Value = 0.7412538312889448
$ws.Range("P12").Value = 0.7412538312889448
$ws.Range("Q12").Value = 2.429437277622778
$ws.Range("R12").Value = 21.864935498605
$ws.Range("S12").Value = 0.0002166157916459741
$ws.Range("T12").Value = 0.0002166157916459741

# Row 13: Inflammatory-Mac -> FAPs
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Col8a1"
$ws.Range("C13").Value = "Itga1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.05623966666666667
$ws.Range("H13").Value = 0.168719
$ws.Range("I13").Value = 0.0002922288998753735
$ws.Range("J13").Value = 0.0002922288998753735
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.38032666666667
$ws.Range("N13").Value = 31.14098
$ws.Range("O13").Value = 0.1781209566020688
$ws.Range("P13").Value = 0.1781209566020688
$ws.Range("Q13").Value = 0.5837861116244445
$ws.Range("R13").Value = 5.25407500462
$ws.Range("S13").Value = 0.00005205209119257172
$ws.Range("T13").Value = 0.00005205209119257172

# Row 14: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("A14").Value = "Inflammatory-Mac"
$ws.Range("B14").Value = "Col8a1"
$ws.Range("C14").Value = "Itga1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05623966666666667
$ws.Range("H14").Value = 0.168719
$ws.Range("I14").Value = 0.0002922288998753735
$ws.Range("J14").Value = 0.0002922288998753735
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.422089
$ws.Range("N14").Value = 1.266267
$ws.Range("O14").Value = 0.007242825670663926
$ws.Range("P14").Value = 0.007242825670663927
$ws.Range("Q14").Value = 0.02373814466366667
$ws.Range("R14").Value = 0.213643301973
$ws.Range("S14").Value = 0.000002116562977727233
$ws.Range("T14").Value = 0.000002116562977727234

# Row 15: Inflammatory-Mac -> MuSCs
$ws.Range("A15").Value = "Inflammatory-Mac"
$ws.Range("B15").Value = "Col8a1"
$ws.Range("C15").Value = "Itga1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05623966666666667
$ws.Range("H15").Value = 0.168719
$ws.Range("I15").Value = 0.0002922288998753735
$ws.Range("J15").Value = 0.0002922288998753735
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.133026333333333
$ws.Range("N15").Value = 12.399079
$ws.Range("O15").Value = 0.07092056230936286
$ws.Range("P15").Value = 0.07092056230936288
$ws.Range("Q15").Value = 0.2324400233112222
$ws.Range("R15").Value = 2.091960209801
$ws.Range("S15").Value = 0.00002072503790220799
$ws.Range("T15").Value = 0.00002072503790220799

# Row 16: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A16").Value = "Inflammatory-Mac"
$ws.Range("B16").Value = "Col8a1"
$ws.Range("C16").Value = "Itga1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05623966666666667
$ws.Range("H16").Value = 0.168719
$ws.Range("I16").Value = 0.0002922288998753735
$ws.Range("J16").Value = 0.0002922288998753735
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.1434673333333333
$ws.Range("N16").Value = 0.430402
$ws.Range("O16").Value = 0.002461824128959449
$ws.Range("P16").Value = 0.002461824128959449
$ws.Range("Q16").Value = 0.008068555004222223
$ws.Range("R16").Value = 0.07261699503800001
$ws.Range("S16").Value = 0.0000007194161568924696
$ws.Range("T16").Value = 0.0000007194161568924696

# Row 17: MuSCs -> ECs
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Col8a1"
$ws.Range("C17").Value = "Itga1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 5.493340666666666
$ws.Range("H17").Value = 16.480022
$ws.Range("I17").Value = 0.02854413965814136
$ws.Range("J17").Value = 0.02854413965814136
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 43.19793166666667
$ws.Range("N17").Value = 129.593795
$ws.Range("O17").Value = 0.7412538312889448
$ws.Range("P17").Value = 0.7412538312889448
$ws.Range("Q17").Value = 237.3009547403877
$ws.Range("R17").Value = 2135.70859266349
$ws.Range("S17").Value = 0.02115845288244399
$ws.Range("T17").Value = 0.02115845288244399

# Row 18: MuSCs -> FAPs
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Col8a1"
$ws.Range("C18").Value = "Itga1"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 5.493340666666666
$ws.Range("H18").Value = 16.480022
$ws.Range("I18").Value = 0.02854413965814136
$ws.Range("J18").Value = 0.02854413965814136
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 10.38032666666667
$ws.Range("N18").Value = 31.14098
$ws.Range("O18").Value = 0.1781209566020688
$ws.Range("P18").Value = 0.1781209566020688
$ws.Range("Q18").Value = 57.02267061128444
$ws.Range("R18").Value = 513.2040355015599
$ws.Range("S18").Value = 0.005084309461291189
$ws.Range("T18").Value = 0.005084309461291189

# Row 19: MuSCs -> Inflammatory-Mac
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Col8a1"
$ws.Range("C19").Value = "Itga1"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 5.493340666666666
$ws.Range("H19").Value = 16.480022
$ws.Range("I19").Value = 0.02854413965814136
$ws.Range("J19").Value = 0.02854413965814136
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.422089
$ws.Range("N19").Value = 1.266267
$ws.Range("O19").Value = 0.007242825670663926
$ws.Range("P19").Value = 0.007242825670663927
$ws.Range("Q19").Value = 2.318678668652666
$ws.Range("R19").Value = 20.868108017874
$ws.Range("S19").Value = 0.0002067402274630024
$ws.Range("T19").Value = 0.0002067402274630025

# Row 20: MuSCs -> MuSCs
$ws.Range("A20").Value = "MuSCs"
$ws.Range("B20").Value = "Col8a1"
$ws.Range("C20").Value = "Itga1"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 5.493340666666666
$ws.Range("H20").Value = 16.480022
$ws.Range("I20").Value = 0.02854413965814136
$ws.Range("J20").Value = 0.02854413965814136
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 4.133026333333333
$ws.Range("N20").Value = 12.399079
$ws.Range("O20").Value = 0.07092056230936286
$ws.Range("P20").Value = 0.07092056230936288
$ws.Range("Q20").Value = 22.70412163330422
$ws.Range("R20").Value = 204.3370946997379
$ws.Range("S20").Value = 0.00202436643519237
$ws.Range("T20").Value = 0.00202436643519237

# Row 21: MuSCs -> Resolving-Mac
$ws.Range("A21").Value = "MuSCs"
$ws.Range("B21").Value = "Col8a1"
$ws.Range("C21").Value = "Itga1"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 5.493340666666666
$ws.Range("H21").Value = 16.480022
$ws.Range("I21").Value = 0.02854413965814136
$ws.Range("J21").Value = 0.02854413965814136
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.1434673333333333
$ws.Range("N21").Value = 0.430402
$ws.Range("O21").Value = 0.002461824128959449
$ws.Range("P21").Value = 0.002461824128959449
$ws.Range("Q21").Value = 0.7881149365382221
$ws.Range("R21").Value = 7.093034428844
$ws.Range("S21").Value = 0.00007027065175080072
$ws.Range("T21").Value = 0.00007027065175080072

# Row 22: Resolving-Mac -> ECs
$ws.Range("A22").Value = "Resolving-Mac"
$ws.Range("B22").Value = "Col8a1"
$ws.Range("C22").Value = "Itga1"
$ws.Range("D22").Value = "ECs"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.05570833333333333
$ws.Range("H22").Value = 0.167125
$ws.Range("I22").Value = 0.0002894680201499048
$ws.Range("J22").Value = 0.0002894680201499048
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 43.19793166666667
$ws.Range("N22").Value = 129.593795
$ws.Range("O22").Value = 0.7412538312889448
$ws.Range("P22").Value = 0.7412538312889448
$ws.Range("Q22").Value = 2.406484776597222
$ws.Range("R22").Value = 21.658362989375
$ws.Range("S22").Value = 0.0002145692789717424
$ws.Range("T22").Value = 0.0002145692789717424

# Row 23: Resolving-Mac -> FAPs
$ws.Range("A23").Value = "Resolving-Mac"
$ws.Range("B23").Value = "Col8a1"
$ws.Range("C23").Value = "Itga1"
$ws.Range("D23").Value = "FAPs"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.05570833333333333
$ws.Range("H23").Value = 0.167125
$ws.Range("I23").Value = 0.0002894680201499048
$ws.Range("J23").Value = 0.0002894680201499048
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 10.38032666666667
$ws.Range("N23").Value = 31.14098
$ws.Range("O23").Value = 0.1781209566020688
$ws.Range("P23").Value = 0.1781209566020688
$ws.Range("Q23").Value = 0.5782706980555555
$ws.Range("R23").Value = 5.2044362825
$ws.Range("S23").Value = 0.00005156032065480798
$ws.Range("T23").Value = 0.00005156032065480798

# Row 24: Resolving-Mac -> Inflammatory-Mac
$ws.Range("A24").Value = "Resolving-Mac"
$ws.Range("B24").Value = "Col8a1"
$ws.Range("C24").Value = "Itga1"
$ws.Range("D24").Value = "Inflammatory-Mac"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.05570833333333333
$ws.Range("H24").Value = 0.167125
$ws.Range("I24").Value = 0.0002894680201499048
$ws.Range("J24").Value = 0.0002894680201499048
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.422089
$ws.Range("N24").Value = 1.266267
$ws.Range("O24").Value = 0.007242825670663926
$ws.Range("P24").Value = 0.007242825670663927
$ws.Range("Q24").Value = 0.02351387470833333
$ws.Range("R24").Value = 0.211624872375
$ws.Range("S24").Value = 0.000002096566407177993
$ws.Range("T24").Value = 0.000002096566407177994

# Row 25: Resolving-Mac -> MuSCs
$ws.Range("A25").Value = "Resolving-Mac"
$ws.Range("B25").Value = "Col8a1"
$ws.Range("C25").Value = "Itga1"
$ws.Range("D25").Value = "MuSCs"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.05570833333333333
$ws.Range("H25").Value = 0.167125
$ws.Range("I25").Value = 0.0002894680201499048
$ws.Range("J25").Value = 0.0002894680201499048
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 4.133026333333333
$ws.Range("N25").Value = 12.399079
$ws.Range("O25").Value = 0.07092056230936286
$ws.Range("P25").Value = 0.07092056230936288
$ws.Range("Q25").Value = 0.2302440086527777
$ws.Range("R25").Value = 2.072196077875
$ws.Range("S25").Value = 0.00002052923475960923
$ws.Range("T25").Value = 0.00002052923475960923

# Row 26: Resolving-Mac -> Resolving-Mac
$ws.Range("A26").Value = "Resolving-Mac"
$ws.Range("B26").Value = "Col8a1"
$ws.Range("C26").Value = "Itga1"
$ws.Range("D26").Value = "Resolving-Mac"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.05570833333333333
$ws.Range("H26").Value = 0.167125
$ws.Range("I26").Value = 0.0002894680201499048
$ws.Range("J26").Value = 0.0002894680201499048
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.1434673333333333
$ws.Range("N26").Value = 0.430402
$ws.Range("O26").Value = 0.002461824128959449
$ws.Range("P26").Value = 0.002461824128959449
$ws.Range("Q26").Value = 0.007992326027777778
$ws.Range("R26").Value = 0.07193093425
$ws.Range("S26").Value = 0.0000007126193565671559
$ws.Range("T26").Value = 0.0000007126193565671559
